$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

function Replace-WholeWord($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "28 août 2023" "dayMonthYear"
Replace-Text "masterName" "lawyerName"
Replace-Text "Madame.l’avocate@avocate.com" "lawyerEmail"
Replace-Text "Madame.l’avocate@L’avocate.com" "lawyerNotificationEmail"
Replace-Text "xxx-xxx-xxxx" "faxNumber"
Replace-Text "450-333-3333 poste 333" "phoneNumber"
Replace-Text "Monsieur le Notifié" "notifieda"
Replace-Text "Courriel" "notifiedEmaila"
Replace-Text "200-04-xxxxxx-xxx" "courtNumber"
Replace-Text "Monsieur Famille c. Madame Famille" "formattedPartsNames"
Replace-Text "Nouvel avis de présentation" "documentName"
Replace-WholeWord "x" "pageCount"
